$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 163, shifting rows 163:237 down to 164:238
$ws.Rows.Item(163).Insert()

# Populate the new row 163 with the new data entry
$ws.Cells.Item(163, 1).Value = 4
$ws.Cells.Item(163, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(163, 3).Value = "Los Lagos"
$ws.Cells.Item(163, 4).Value = 44609
$ws.Cells.Item(163, 5).Value = 10
$ws.Cells.Item(163, 6).Value = 100112037
$ws.Cells.Item(163, 7).Value = "Cebollín"
$ws.Cells.Item(163, 8).Value = "Sin especificar"
$ws.Cells.Item(163, 9).Value = "Primera"
$ws.Cells.Item(163, 10).Value = 60
$ws.Cells.Item(163, 11).Value = 6000
$ws.Cells.Item(163, 12).Value = 6500
$ws.Cells.Item(163, 13).Value = 6250
$ws.Cells.Item(163, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(163, 15).Value = "Región Metropolitana"
$ws.Cells.Item(163, 16).Value = 174
$ws.Cells.Item(163, 17).Value = 36
$ws.Cells.Item(163, 18).Value = "Hortaliza"
